$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top for a header label
$ws.Rows.Item(1).Insert()

# Set the new header text
$ws.Range("A1").Value = "Course outcome"

# Bold the new header cell
$ws.Range("A1").Font.Bold = $true

# Move the active selection (matches final cursor position after editing)
[void]$ws.Range("E6").Select()

